# Adjust character physic material, adjust attackbalance
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arkusz1")

# --- Row 3 (slp): F3 30 -> 29 (C3 recalculates from 13 -> 12) ---
$ws.Range("F3").Value = 29

# --- Row 5 (srk): F5 24 -> 26, G5 15 -> 16 (C5 -1 -> 1, D5 -15 -> -14) ---
$ws.Range("F5").Value = 26
$ws.Range("G5").Value = 16

# --- Row 6 (clp): F6 37 -> 33 (C6 16 -> 12) ---
$ws.Range("F6").Value = 33

# --- Row 7 (crp): F7 12 -> 18 (C7 -3 -> 3) ---
$ws.Range("F7").Value = 18

# --- Row 25 (srk, mirrored table): C25 -1 -> 1, D25 -15 -> -14 (F25/G25 recalc 24/15 -> 26/16) ---
$ws.Range("C25").Value = 1
$ws.Range("D25").Value = -14

# --- Row 26 (clp, mirrored table): C26 16 -> 12 (F26 recalcs 37 -> 33) ---
$ws.Range("C26").Value = 12

# --- Row 27 (crp, mirrored table): C27 -3 -> 3 (F27 recalcs 12 -> 18) ---
$ws.Range("C27").Value = 3

# --- View state: scroll so row 10 is the top-left visible row, select C23 ---
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C23").Select()
